# Insert a new data row before the current row 40. This shifts the
# existing rows 40-108 down to 41-109 (preserving all of their values
# and the date-column style), and leaves a blank row 40 for the new
# weekly price observation to be filled in below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("40:40").Insert()

$ws.Range("A40").Value = 6
$ws.Range("B40").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C40").Value = "Metropolitana"
$ws.Range("D40").Value = 44935
$ws.Range("E40").Value = 13
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100101
$ws.Range("H40").Value = "Berries"
$ws.Range("I40").Value = 100101008
$ws.Range("J40").Value = "Mora"
$ws.Range("K40").Value = "Sin especificar"
$ws.Range("L40").Value = "Primera"
$ws.Range("M40").Value = 100
$ws.Range("N40").Value = 4000
$ws.Range("O40").Value = 4000
$ws.Range("P40").Value = 4000
$ws.Range("Q40").Value = "`$/bandeja 2 kilos"
$ws.Range("R40").Value = "Región del Maule"
$ws.Range("S40").Value = 2000
$ws.Range("T40").Value = 2
